# "Fruta / hortaliza, semanal" - weekly refresh of the Acelga / Vega Modelo
# de Temuco consolidated price sheet: a new daily record is inserted at
# row 123, pushing the previously-existing rows 123:189 down to 124:190
# (dimension grows from A1:R189 to A1:R190). No other existing data changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 123; Excel shifts rows 123:189 -> 124:190 and
# carries the row-123 formatting (e.g. the date style on column D) down
# with the cells that move, just like a normal Excel "Insert Row".
$ws.Rows.Item(123).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A123").Value = 10
$ws.Range("B123").Value = "Vega Modelo de Temuco"
$ws.Range("C123").Value = "La Araucanía"
$ws.Range("D123").Value = 44460
$ws.Range("E123").Value = 9
$ws.Range("F123").Value = 100112009
$ws.Range("G123").Value = "Acelga"
$ws.Range("H123").Value = "Sin especificar"
$ws.Range("I123").Value = "Primera"
$ws.Range("J123").Value = 20
$ws.Range("K123").Value = 8000
$ws.Range("L123").Value = 8000
$ws.Range("M123").Value = 8000
$ws.Range("N123").Value = "$/docena de atados (12 kilos)"
$ws.Range("O123").Value = "Provincia de Cautín"
$ws.Range("P123").Value = 667
$ws.Range("Q123").Value = 12
$ws.Range("R123").Value = "Hortaliza"
